# Default duration of class is 1 (hour) instead of 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance")

$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D4").Value = 1
